# The edit removes the "实际结果" (F) and "是否通过" (G) data columns'
# contents for the test-case rows (2-9). The "预期结果" (E) column values
# are left untouched. Column G's first data cell (G2) keeps its existing
# cell style but becomes blank; the remaining G cells (G3:G9) and all of
# column F's data cells (F2:F9) are cleared entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "实际结果" (actual result) values in column F, rows 2-9.
$ws.Range("F2:F9").ClearContents()

# Remove the "是否通过" (pass/fail) values in column G, rows 2-9.
$ws.Range("G2:G9").ClearContents()

$wb.Save()
